# Add a new "Quantity" (數量) column (F) to the order-tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 --------------------------------------------------
# Re-use the same font ("微軟正黑體") already applied to the other
# header cells (A1/B1/C1/E1) so the style table gains no extra font and
# the new header xf matches the one used elsewhere in the sheet.
$headerFontName = $ws.Range("A1").Font.Name

$ws.Range("F1").Value = "數量"
$ws.Range("F1").Font.Name = $headerFontName

# --- Data cells F2:F21 (quantity per order row) -----------------------
$quantities = @(1, 11, 3, 4, 11, 6, 7, 8, 20, 10, 11, 7, 13, 14, 15, 20, 17, 7, 19, 20)

for ($i = 0; $i -lt $quantities.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $quantities[$i]
}

# --- Selection ----------------------------------------------------------
$ws.Range("I12").Select()
